$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C13: 1.0 -> 0.0
$ws.Range("C13").Value = 0.0

# D13: empty -> "1.5.3"
$ws.Range("D13").Value = "1.5.3"
